$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Mes" date column (C2:C41) from 1-Jun-2023 (45078) to 1-Jul-2023 (45108)
$ws.Range("C2:C41").Value2 = 45108

# 2. Update the "Capacidad instalada" column (D) values that changed
$dChanges = @{
    2  = 516
    4  = 360
    9  = 80
    10 = 72
    11 = 216
    14 = 328
    15 = 328
    16 = 141
    17 = 2056
    18 = 177
    20 = 344
    23 = 106
    24 = 516
    25 = 70
    26 = 48
    27 = 300
    28 = 369
    29 = 0
    30 = 328
    32 = 21
    33 = 516
}

foreach ($row in $dChanges.Keys) {
    $ws.Cells.Item($row, 4).Value = $dChanges[$row]
}

# 3. Remove the trailing rows (42-57) that held extra months of data no longer needed
$ws.Range("A42:A57").EntireRow.Delete()

# 4. Reset the view: clear the scrolled/frozen top-left cell and move the selection
$ws.Range("D42").Select()
